# init templates for expense fragment
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10, pushing the SPRINT 4 stories down one row.
$ws.Rows.Item(10).Insert()

# Fill in the new row 10: SPRINT 3 / new story / NOT STARTED
$ws.Range("A10").Value = "SPRINT 3"
$ws.Range("B10").Value = "I want to copy over existing expenses to new month"
$ws.Range("C10").Value = "NOT STARTED"

# Copy the look (styles) of the row above (row 9) onto the new row 10
$ws.Range("A9:C9").Copy()
$ws.Range("A10:C10").PasteSpecial(-4122) # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Re-order the two stories that followed SPRINT 4's first item so "validate input"
# moves up to directly follow SPRINT 3, and the avatar/combined-image stories shift down.
$ws.Range("B11").Value = "I want to validate input when adding new room"
$ws.Range("B12").Value = "I want to get avatar for each housemate from the phone contact"
$ws.Range("B13").Value = "I want to dispay room avatar as a combined image from roommates' avatars"

# Highlight the in-progress sprint group (rows 5-7, column A) with the same
# "Neutral" look already used by C5 (fill + border), rather than the
# Style-assignment shortcut (which creates a borderless variant).
$ws.Range("C5").Copy()
$ws.Range("A5:A7").PasteSpecial(-4122) # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Update the table range to include the new row
$table = $ws.ListObjects.Item("Table2")
$table.Resize($ws.Range("A1:C13"))

# Update selection to match the authored state
$ws.Range("B11").Select()
